$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H64").Value = 58422.223
$ws.Range("J64").Value = 3050
$ws.Range("L64").Value = 3050
$ws.Range("N64").Value = -3546

$ws.Range("H67").Value = 58422.223
$ws.Range("J67").Value = 3050
$ws.Range("L67").Value = 3050
$ws.Range("N67").Value = -4766

$ws.Range("H129").Value = 2398.6287
$ws.Range("J129").Value = 1046.9482
$ws.Range("L129").Value = 3140.8446
$ws.Range("N129").Value = -13140.8446

$ws.Range("H132").Value = 3850352.8
$ws.Range("I132").Value = 4549780.5
$ws.Range("J132").Value = 3499.6
$ws.Range("K132").Value = 13649341.5
$ws.Range("L132").Value = 10498.8
$ws.Range("M132").Value = -13646811.5
$ws.Range("N132").Value = -15558.8

$ws.Range("H135").Value = 1573.9556
$ws.Range("I135").Value = 707.129
$ws.Range("J135").Value = 3493.3572
$ws.Range("K135").Value = 6364.161
$ws.Range("L135").Value = 31440.2148
$ws.Range("M135").Value = -3829.161
$ws.Range("N135").Value = -36510.2148

$ws.Range("H137").Value = 1525.2106
$ws.Range("I137").Value = 1241.6552
$ws.Range("K137").Value = 3724.9656
$ws.Range("M137").Value = -1174.9656

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 24448.512
$ws.Range("I32").Value = 6719.3
$ws.Range("J32").Value = 260838
$ws.Range("K32").Value = 6719.3
$ws.Range("L32").Value = 260838
$ws.Range("M32").Value = -6432.3
$ws.Range("N32").Value = -261412

$ws.Range("H45").Value = 202532
$ws.Range("I45").Value = 334233.34
$ws.Range("J45").Value = 4980
$ws.Range("K45").Value = 334233.34
$ws.Range("L45").Value = 4980
$ws.Range("M45").Value = -333856.34
$ws.Range("N45").Value = -5734

$ws.Range("H74").Value = 1065.8667
$ws.Range("I74").Value = 997.8182
$ws.Range("K74").Value = 997.8182
$ws.Range("M74").Value = -123.8182

$ws.Range("H77").Value = 1065.8667
$ws.Range("I77").Value = 997.8182
$ws.Range("K77").Value = 4989.091
$ws.Range("M77").Value = -621.0910000000003

$ws.Range("H88").Value = 1000
$ws.Range("J88").Value = 1000
$ws.Range("L88").Value = 1000
$ws.Range("N88").Value = -1812

$ws.Range("H91").Value = 1000
$ws.Range("J91").Value = 1000
$ws.Range("L91").Value = 1000
$ws.Range("N91").Value = -3808

$ws.Range("H122").Value = 2725.7778
$ws.Range("I122").Value = 2186
$ws.Range("J122").Value = 4615
$ws.Range("K122").Value = 6558
$ws.Range("L122").Value = 13845
$ws.Range("M122").Value = -4108
$ws.Range("N122").Value = -18745

$ws.Range("H132").Value = 13014.115
$ws.Range("I132").Value = 15943.925
$ws.Range("K132").Value = 47831.77499999999
$ws.Range("M132").Value = -45301.77499999999

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 51750.273
$ws.Range("I86").Value = 62778.11
$ws.Range("J86").Value = 2125
$ws.Range("K86").Value = 62778.11
$ws.Range("L86").Value = 2125
$ws.Range("M86").Value = -61655.11
$ws.Range("N86").Value = -4371

$ws.Range("H89").Value = 51750.273
$ws.Range("I89").Value = 62778.11
$ws.Range("J89").Value = 2125
$ws.Range("K89").Value = 313890.55
$ws.Range("L89").Value = 10625
$ws.Range("M89").Value = -308274.55
$ws.Range("N89").Value = -21857

$ws.Range("H134").Value = 2656.4792
$ws.Range("I134").Value = 2329.818
$ws.Range("K134").Value = 6989.454000000001
$ws.Range("M134").Value = -4454.454000000001

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 84260.25
$ws.Range("I16").Value = 1001.4286
$ws.Range("J16").Value = 200822.6
$ws.Range("K16").Value = 1001.4286
$ws.Range("L16").Value = 200822.6
$ws.Range("M16").Value = -714.4286
$ws.Range("N16").Value = -201396.6

$ws.Range("H31").Value = 36388.832
$ws.Range("I31").Value = 807.9231
$ws.Range("J31").Value = 52338.9
$ws.Range("K31").Value = 807.9231
$ws.Range("L31").Value = 52338.9
$ws.Range("M31").Value = -512.9231
$ws.Range("N31").Value = -52928.9

$ws.Range("H34").Value = 36388.832
$ws.Range("I34").Value = 807.9231
$ws.Range("J34").Value = 52338.9
$ws.Range("K34").Value = 807.9231
$ws.Range("L34").Value = 52338.9
$ws.Range("M34").Value = -605.9231
$ws.Range("N34").Value = -52742.9

$ws.Range("H113").Value = 84260.25
$ws.Range("I113").Value = 1001.4286
$ws.Range("J113").Value = 200822.6
$ws.Range("K113").Value = 1001.4286
$ws.Range("L113").Value = 200822.6
$ws.Range("M113").Value = 1168.5714
$ws.Range("N113").Value = -205162.6

$ws.Range("H132").Value = 4609.857
$ws.Range("I132").Value = 5191.778
$ws.Range("J132").Value = 3562.4
$ws.Range("K132").Value = 15575.334
$ws.Range("L132").Value = 10687.2
$ws.Range("M132").Value = -13045.334
$ws.Range("N132").Value = -15747.2

$ws.Range("H134").Value = 1368.909
$ws.Range("I134").Value = 850.8889
$ws.Range("J134").Value = 3700
$ws.Range("K134").Value = 2552.6667
$ws.Range("L134").Value = 11100
$ws.Range("M134").Value = -17.66670000000022
$ws.Range("N134").Value = -16170

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 783.9231
$ws.Range("I2").Value = 17.4
$ws.Range("J2").Value = 1263
$ws.Range("K2").Value = 104.4
$ws.Range("L2").Value = 7578
$ws.Range("M2").Value = 8.600000000000009
$ws.Range("N2").Value = -7804

$ws.Range("H98").Value = 63439.438
$ws.Range("J98").Value = 67615.39999999999
$ws.Range("L98").Value = 202846.2
$ws.Range("N98").Value = -205842.2

$ws.Range("H129").Value = 5880.7393
$ws.Range("I129").Value = 595.8
$ws.Range("J129").Value = 7348.778
$ws.Range("K129").Value = 1787.4
$ws.Range("L129").Value = 22046.334
$ws.Range("M129").Value = 3212.6
$ws.Range("N129").Value = -32046.334

$ws.Range("H131").Value = 866449.3
$ws.Range("I131").Value = 603.9
$ws.Range("J131").Value = 1043152.44
$ws.Range("K131").Value = 1811.7
$ws.Range("L131").Value = 3129457.32
$ws.Range("M131").Value = 3228.3
$ws.Range("N131").Value = -3139537.32

$ws.Range("H138").Value = 7749.6875
$ws.Range("I138").Value = 9666.416999999999
$ws.Range("J138").Value = 1999.5
$ws.Range("K138").Value = 28999.251
$ws.Range("L138").Value = 5998.5
$ws.Range("M138").Value = -23859.251
$ws.Range("N138").Value = -16278.5

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 58311.133
$ws.Range("I70").Value = 85207.88
$ws.Range("K70").Value = 85207.88
$ws.Range("M70").Value = -84937.88

$ws.Range("H73").Value = 58311.133
$ws.Range("I73").Value = 85207.88
$ws.Range("K73").Value = 85207.88
$ws.Range("M73").Value = -84271.88

$ws.Range("H122").Value = 3243.7
$ws.Range("I122").Value = 2553.625
$ws.Range("K122").Value = 7660.875
$ws.Range("M122").Value = -5210.875

$ws.Range("H123").Value = 30000
$ws.Range("J123").Value = 30000
$ws.Range("L123").Value = 30000
$ws.Range("N123").Value = -34900

$ws.Range("H132").Value = 3511.3547
$ws.Range("I132").Value = 2632.55
$ws.Range("K132").Value = 7897.650000000001
$ws.Range("M132").Value = -5367.650000000001

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H68").Value = 3083.2
$ws.Range("I68").Value = 2052.261
$ws.Range("J68").Value = 5059.1665
$ws.Range("K68").Value = 2052.261
$ws.Range("L68").Value = 5059.1665
$ws.Range("M68").Value = -1303.261
$ws.Range("N68").Value = -6557.1665

$ws.Range("H71").Value = 3083.2
$ws.Range("I71").Value = 2052.261
$ws.Range("J71").Value = 5059.1665
$ws.Range("K71").Value = 10261.305
$ws.Range("L71").Value = 25295.8325
$ws.Range("M71").Value = -6517.305
$ws.Range("N71").Value = -32783.8325

$ws.Range("H132").Value = 3282.4583
$ws.Range("I132").Value = 3282.4583
$ws.Range("K132").Value = 9847.374899999999
$ws.Range("M132").Value = -7317.374899999999

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 6763.625
$ws.Range("I132").Value = 4138.2856
$ws.Range("K132").Value = 12414.8568
$ws.Range("M132").Value = -9884.856800000001
